# Texturometer forces: add experimental data collected on 230718.
# Each replicate (FF1..FF5) has sub-samples A..F; F20/F80 force readings are
# recorded as text (as in the rest of the sheet) in columns C/D.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -- Pass 1: write the new sample Ids (column A) for the "A/B/C" sub-samples first,
#    matching the order the source data was collected in.
$ws.Range("A143").Value = "230718_FF1_A"
$ws.Range("A144").Value = "230718_FF1_B"
$ws.Range("A145").Value = "230718_FF1_C"
$ws.Range("A148").Value = "230718_FF2_A"
$ws.Range("A149").Value = "230718_FF2_B"
$ws.Range("A150").Value = "230718_FF2_C"
$ws.Range("A153").Value = "230718_FF3_A"
$ws.Range("A154").Value = "230718_FF3_B"
$ws.Range("A155").Value = "230718_FF3_C"
$ws.Range("A157").Value = "230718_FF4_A"
$ws.Range("A158").Value = "230718_FF4_B"
$ws.Range("A159").Value = "230718_FF4_C"
$ws.Range("A161").Value = "230718_FF5_A"
$ws.Range("A162").Value = "230718_FF5_B"
$ws.Range("A163").Value = "230718_FF5_C"

# -- Pass 2: write the sample Ids (column A) for the "D/E/F" sub-samples.
$ws.Range("A146").Value = "230718_FF1_D"
$ws.Range("A147").Value = "230718_FF1_E"
$ws.Range("A151").Value = "230718_FF2_D"
$ws.Range("A152").Value = "230718_FF2_E"
$ws.Range("A156").Value = "230718_FF3_D"
$ws.Range("A160").Value = "230718_FF4_D"
$ws.Range("A164").Value = "230718_FF5_D"
$ws.Range("A165").Value = "230718_FF5_E"
$ws.Range("A166").Value = "230718_FF5_F"

# -- Pass 3: fill in date (B), F20/F80 force readings (C/D, kept as text like the
#    rest of the column) and the failed flag (E) for every new row, top to bottom.
$ws.Range("B143").Value = 230718
$ws.Range("C143").Value = "'4.06632804870642"
$ws.Range("D143").Value = "'97.6550598144529"
$ws.Range("E143").Value = 0
$ws.Range("B144").Value = 230718
$ws.Range("C144").Value = "'16.2561569213866"
$ws.Range("D144").Value = "'50.7060813903802"
$ws.Range("E144").Value = 0
$ws.Range("B145").Value = 230718
$ws.Range("C145").Value = "'15.254155158996"
$ws.Range("D145").Value = "'28.4071407318117"
$ws.Range("E145").Value = 0
$ws.Range("B146").Value = 230718
$ws.Range("C146").Value = "'7.65582656860374"
$ws.Range("D146").Value = "'53.2443695068363"
$ws.Range("E146").Value = 0
$ws.Range("B147").Value = 230718
$ws.Range("C147").Value = "'19.1822929382315"
$ws.Range("D147").Value = "'61.8381805419913"
$ws.Range("E147").Value = 0
$ws.Range("B148").Value = 230718
$ws.Range("C148").Value = "'8.89438629150381"
$ws.Range("D148").Value = "'45.6641693115232"
$ws.Range("E148").Value = 0
$ws.Range("B149").Value = 230718
$ws.Range("C149").Value = "'4.27507400512669"
$ws.Range("D149").Value = "'47.7472496032698"
$ws.Range("E149").Value = 0
$ws.Range("B150").Value = 230718
$ws.Range("C150").Value = "'19.6130981445324"
$ws.Range("D150").Value = "'49.5825653076173"
$ws.Range("E150").Value = 0
$ws.Range("B151").Value = 230718
$ws.Range("C151").Value = "'7.58116197585891"
$ws.Range("D151").Value = "'48.5847244262693"
$ws.Range("E151").Value = 0
$ws.Range("B152").Value = 230718
$ws.Range("C152").Value = "'20.1609725952161"
$ws.Range("D152").Value = "'74.8009185791046"
$ws.Range("E152").Value = 0
$ws.Range("B153").Value = 230718
$ws.Range("C153").Value = "'4.770984172821"
$ws.Range("D153").Value = "'72.1743392944364"
$ws.Range("E153").Value = 0
$ws.Range("B154").Value = 230718
$ws.Range("C154").Value = "'9.43238544464287"
$ws.Range("D154").Value = "'55.0407867431644"
$ws.Range("E154").Value = 0
$ws.Range("B155").Value = 230718
$ws.Range("C155").Value = "'11.138648033142"
$ws.Range("D155").Value = "'101.482803344719"
$ws.Range("E155").Value = 0
$ws.Range("B156").Value = 230718
$ws.Range("C156").Value = "'26.9608497619636"
$ws.Range("D156").Value = "'67.6510696411109"
$ws.Range("E156").Value = 0
$ws.Range("B157").Value = 230718
$ws.Range("C157").Value = "'13.2886114120482"
$ws.Range("D157").Value = "'81.6852111816406"
$ws.Range("E157").Value = 0
$ws.Range("B158").Value = 230718
$ws.Range("C158").Value = "'14.3018760681141"
$ws.Range("D158").Value = "'49.6651878356928"
$ws.Range("E158").Value = 0
$ws.Range("B159").Value = 230718
$ws.Range("C159").Value = "'2.9882698059083"
$ws.Range("D159").Value = "'24.8318862915028"
$ws.Range("E159").Value = 0
$ws.Range("B160").Value = 230718
$ws.Range("C160").Value = "'14.8666410446173"
$ws.Range("D160").Value = "'33.0632324218743"
$ws.Range("E160").Value = 0
$ws.Range("B161").Value = 230718
$ws.Range("C161").Value = "'16.0089015960688"
$ws.Range("D161").Value = "'43.2277221679739"
$ws.Range("E161").Value = 0
$ws.Range("B162").Value = 230718
$ws.Range("C162").Value = "'25.1563796997044"
$ws.Range("D162").Value = "'48.6795845031718"
$ws.Range("E162").Value = 0
$ws.Range("B163").Value = 230718
$ws.Range("C163").Value = "'10.8912563323974"
$ws.Range("D163").Value = "'27.513494491578"
$ws.Range("E163").Value = 0
$ws.Range("B164").Value = 230718
$ws.Range("C164").Value = "'6.42561435699454"
$ws.Range("D164").Value = "'45.5603981018074"
$ws.Range("E164").Value = 0
$ws.Range("B165").Value = 230718
$ws.Range("C165").Value = "'16.1440563201911"
$ws.Range("D165").Value = "'47.3149795532241"
$ws.Range("E165").Value = 0
$ws.Range("B166").Value = 230718
$ws.Range("C166").Value = "'11.783949851989"
$ws.Range("D166").Value = "'30.7021770477293"
$ws.Range("E166").Value = 0

# -- Column A is now a little wider ("230718_FFx_y" ids are longer); resize it.
$ws.Columns("A").ColumnWidth = 14.85

# -- Leave the selection/scroll position where the user ended up after entering the data.
$ws.Range("G163").Select()
